$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 23.77588366666667
$ws.Cells.Item(2, 8).Value = 71.327651
$ws.Cells.Item(2, 9).Value = 0.201093431146956
$ws.Cells.Item(2, 10).Value = 0.2010934311469559
$ws.Cells.Item(2, 13).Value = 45.90594266666667
$ws.Cells.Item(2, 14).Value = 137.717828
$ws.Cells.Item(2, 15).Value = 0.3954672001633582
$ws.Cells.Item(2, 16).Value = 0.3954672001633583
$ws.Cells.Item(2, 17).Value = 1091.454352451336
$ws.Cells.Item(2, 18).Value = 9823.089172062028
$ws.Cells.Item(2, 19).Value = 0.07952585618692971
$ws.Cells.Item(2, 20).Value = 0.07952585618692973

# Row 3
$ws.Cells.Item(3, 7).Value = 23.77588366666667
$ws.Cells.Item(3, 8).Value = 71.327651
$ws.Cells.Item(3, 9).Value = 0.201093431146956
$ws.Cells.Item(3, 10).Value = 0.2010934311469559
$ws.Cells.Item(3, 15).Value = 0.3484294080560655
$ws.Cells.Item(3, 16).Value = 0.3484294080560656
$ws.Cells.Item(3, 17).Value = 961.6342234899498
$ws.Cells.Item(3, 18).Value = 8654.708011409548
$ws.Cells.Item(3, 19).Value = 0.07006686517849704
$ws.Cells.Item(3, 20).Value = 0.07006686517849704

# Row 4
$ws.Cells.Item(4, 7).Value = 23.77588366666667
$ws.Cells.Item(4, 8).Value = 71.327651
$ws.Cells.Item(4, 9).Value = 0.201093431146956
$ws.Cells.Item(4, 10).Value = 0.2010934311469559
$ws.Cells.Item(4, 13).Value = 12.761795
$ws.Cells.Item(4, 14).Value = 38.28538500000001
$ws.Cells.Item(4, 15).Value = 0.1099393900775594
$ws.Cells.Item(4, 16).Value = 0.1099393900775594
$ws.Cells.Item(4, 17).Value = 303.4229532978484
$ws.Cells.Item(4, 18).Value = 2730.806579680635
$ws.Cells.Item(4, 19).Value = 0.02210808916890003
$ws.Cells.Item(4, 20).Value = 0.02210808916890003

# Row 5
$ws.Cells.Item(5, 7).Value = 23.77588366666667
$ws.Cells.Item(5, 8).Value = 71.327651
$ws.Cells.Item(5, 9).Value = 0.201093431146956
$ws.Cells.Item(5, 10).Value = 0.2010934311469559
$ws.Cells.Item(5, 13).Value = 16.966758
$ws.Cells.Item(5, 14).Value = 50.900274
$ws.Cells.Item(5, 15).Value = 0.1461640017030168
$ws.Cells.Item(5, 16).Value = 0.1461640017030168
$ws.Cells.Item(5, 17).Value = 403.399664408486
$ws.Cells.Item(5, 18).Value = 3630.596979676374
$ws.Cells.Item(5, 19).Value = 0.02939262061262916
$ws.Cells.Item(5, 20).Value = 0.02939262061262916

# Row 6
$ws.Cells.Item(6, 9).Value = 0.2703947904457373
$ws.Cells.Item(6, 10).Value = 0.2703947904457373
$ws.Cells.Item(6, 13).Value = 45.90594266666667
$ws.Cells.Item(6, 14).Value = 137.717828
$ws.Cells.Item(6, 15).Value = 0.3954672001633582
$ws.Cells.Item(6, 16).Value = 0.3954672001633583
$ws.Cells.Item(6, 17).Value = 1467.594288032687
$ws.Cells.Item(6, 18).Value = 13208.34859229419
$ws.Cells.Item(6, 19).Value = 0.1069322707163337
$ws.Cells.Item(6, 20).Value = 0.1069322707163337

# Row 7
$ws.Cells.Item(7, 9).Value = 0.2703947904457373
$ws.Cells.Item(7, 10).Value = 0.2703947904457373
$ws.Cells.Item(7, 15).Value = 0.3484294080560655
$ws.Cells.Item(7, 16).Value = 0.3484294080560656
$ws.Cells.Item(7, 18).Value = 11637.31677244074
$ws.Cells.Item(7, 19).Value = 0.09421349677645215
$ws.Cells.Item(7, 20).Value = 0.09421349677645216

# Row 8
$ws.Cells.Item(8, 9).Value = 0.2703947904457373
$ws.Cells.Item(8, 10).Value = 0.2703947904457373
$ws.Cells.Item(8, 13).Value = 12.761795
$ws.Cells.Item(8, 14).Value = 38.28538500000001
$ws.Cells.Item(8, 15).Value = 0.1099393900775594
$ws.Cells.Item(8, 16).Value = 0.1099393900775594
$ws.Cells.Item(8, 17).Value = 407.9893878455034
$ws.Cells.Item(8, 18).Value = 3671.904490609531
$ws.Cells.Item(8, 19).Value = 0.02972703834175385
$ws.Cells.Item(8, 20).Value = 0.02972703834175386

# Row 9
$ws.Cells.Item(9, 9).Value = 0.2703947904457373
$ws.Cells.Item(9, 10).Value = 0.2703947904457373
$ws.Cells.Item(9, 13).Value = 16.966758
$ws.Cells.Item(9, 14).Value = 50.900274
$ws.Cells.Item(9, 15).Value = 0.1461640017030168
$ws.Cells.Item(9, 16).Value = 0.1461640017030168
$ws.Cells.Item(9, 17).Value = 542.420342133908
$ws.Cells.Item(9, 18).Value = 4881.783079205172
$ws.Cells.Item(9, 19).Value = 0.03952198461119762
$ws.Cells.Item(9, 20).Value = 0.03952198461119763

# Row 10
$ws.Cells.Item(10, 7).Value = 14.51831366666667
$ws.Cells.Item(10, 8).Value = 43.554941
$ws.Cells.Item(10, 9).Value = 0.1227940694288843
$ws.Cells.Item(10, 10).Value = 0.1227940694288843
$ws.Cells.Item(10, 13).Value = 45.90594266666667
$ws.Cells.Item(10, 14).Value = 137.717828
$ws.Cells.Item(10, 15).Value = 0.3954672001633582
$ws.Cells.Item(10, 16).Value = 0.3954672001633583
$ws.Cells.Item(10, 17).Value = 666.4768747986831
$ws.Cells.Item(10, 18).Value = 5998.291873188148
$ws.Cells.Item(10, 19).Value = 0.04856102683370589
$ws.Cells.Item(10, 20).Value = 0.0485610268337059

# Row 11
$ws.Cells.Item(11, 7).Value = 14.51831366666667
$ws.Cells.Item(11, 8).Value = 43.554941
$ws.Cells.Item(11, 9).Value = 0.1227940694288843
$ws.Cells.Item(11, 10).Value = 0.1227940694288843
$ws.Cells.Item(11, 15).Value = 0.3484294080560655
$ws.Cells.Item(11, 16).Value = 0.3484294080560656
$ws.Cells.Item(11, 17).Value = 587.2045592484964
$ws.Cells.Item(11, 18).Value = 5284.841033236467
$ws.Cells.Item(11, 19).Value = 0.04278506492390157
$ws.Cells.Item(11, 20).Value = 0.04278506492390158

# Row 12
$ws.Cells.Item(12, 7).Value = 14.51831366666667
$ws.Cells.Item(12, 8).Value = 43.554941
$ws.Cells.Item(12, 9).Value = 0.1227940694288843
$ws.Cells.Item(12, 10).Value = 0.1227940694288843
$ws.Cells.Item(12, 13).Value = 12.761795
$ws.Cells.Item(12, 14).Value = 38.28538500000001
$ws.Cells.Item(12, 15).Value = 0.1099393900775594
$ws.Cells.Item(12, 16).Value = 0.1099393900775594
$ws.Cells.Item(12, 17).Value = 185.2797427596983
$ws.Cells.Item(12, 18).Value = 1667.517684837285
$ws.Cells.Item(12, 19).Value = 0.01349990509815302
$ws.Cells.Item(12, 20).Value = 0.01349990509815303

# Row 13
$ws.Cells.Item(13, 7).Value = 14.51831366666667
$ws.Cells.Item(13, 8).Value = 43.554941
$ws.Cells.Item(13, 9).Value = 0.1227940694288843
$ws.Cells.Item(13, 10).Value = 0.1227940694288843
$ws.Cells.Item(13, 13).Value = 16.966758
$ws.Cells.Item(13, 14).Value = 50.900274
$ws.Cells.Item(13, 15).Value = 0.1461640017030168
$ws.Cells.Item(13, 16).Value = 0.1461640017030168
$ws.Cells.Item(13, 17).Value = 246.328714550426
$ws.Cells.Item(13, 18).Value = 2216.958430953834
$ws.Cells.Item(13, 19).Value = 0.01794807257312381
$ws.Cells.Item(13, 20).Value = 0.01794807257312381

# Row 14
$ws.Cells.Item(14, 7).Value = 47.96923
$ws.Cells.Item(14, 8).Value = 143.90769
$ws.Cells.Item(14, 9).Value = 0.4057177089784224
$ws.Cells.Item(14, 10).Value = 0.4057177089784224
$ws.Cells.Item(14, 13).Value = 45.90594266666667
$ws.Cells.Item(14, 14).Value = 137.717828
$ws.Cells.Item(14, 15).Value = 0.3954672001633582
$ws.Cells.Item(14, 16).Value = 0.3954672001633583
$ws.Cells.Item(14, 17).Value = 2202.072722144147
$ws.Cells.Item(14, 18).Value = 19818.65449929732
$ws.Cells.Item(14, 19).Value = 0.1604480464263889
$ws.Cells.Item(14, 20).Value = 0.1604480464263889

# Row 15
$ws.Cells.Item(15, 7).Value = 47.96923
$ws.Cells.Item(15, 8).Value = 143.90769
$ws.Cells.Item(15, 9).Value = 0.4057177089784224
$ws.Cells.Item(15, 10).Value = 0.4057177089784224
$ws.Cells.Item(15, 15).Value = 0.3484294080560655
$ws.Cells.Item(15, 16).Value = 0.3484294080560656
$ws.Cells.Item(15, 17).Value = 1940.153051267347
$ws.Cells.Item(15, 18).Value = 17461.37746140612
$ws.Cells.Item(15, 19).Value = 0.1413639811772148
$ws.Cells.Item(15, 20).Value = 0.1413639811772148

# Row 16
$ws.Cells.Item(16, 7).Value = 47.96923
$ws.Cells.Item(16, 8).Value = 143.90769
$ws.Cells.Item(16, 9).Value = 0.4057177089784224
$ws.Cells.Item(16, 10).Value = 0.4057177089784224
$ws.Cells.Item(16, 13).Value = 12.761795
$ws.Cells.Item(16, 14).Value = 38.28538500000001
$ws.Cells.Item(16, 15).Value = 0.1099393900775594
$ws.Cells.Item(16, 16).Value = 0.1099393900775594
$ws.Cells.Item(16, 17).Value = 612.1734795678501
$ws.Cells.Item(16, 18).Value = 5509.561316110651
$ws.Cells.Item(16, 19).Value = 0.04460435746875251
$ws.Cells.Item(16, 20).Value = 0.04460435746875251

# Row 17
$ws.Cells.Item(17, 7).Value = 47.96923
$ws.Cells.Item(17, 8).Value = 143.90769
$ws.Cells.Item(17, 9).Value = 0.4057177089784224
$ws.Cells.Item(17, 10).Value = 0.4057177089784224
$ws.Cells.Item(17, 13).Value = 16.966758
$ws.Cells.Item(17, 14).Value = 50.900274
$ws.Cells.Item(17, 15).Value = 0.1461640017030168
$ws.Cells.Item(17, 16).Value = 0.1461640017030168
$ws.Cells.Item(17, 17).Value = 813.88231685634
$ws.Cells.Item(17, 18).Value = 7324.94085170706
$ws.Cells.Item(17, 19).Value = 0.0593013239060662
$ws.Cells.Item(17, 20).Value = 0.05930132390606622

